$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025".
#    Setting .Text on the *whole* paragraph Range (the one that already owns
#    a run) keeps the xml:space="preserve" attribute on the <w:t> element,
#    matching how the source document already had it.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*September 19, 2025*") {
        $p.Range.Text = "September 21, 2025"
        break
    }
}

# ---------------------------------------------------------------------------
# 2. Split the single-line mailing address paragraph into two paragraphs:
#      "2900 Sanor Pl"
#      "Santa Clara, CA 95051"
#    (Only the letter-header occurrence is touched, not the identical text
#    that also appears inside the property-details table further down.)
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*2900 Sanor Pl, Santa Clara CA 95051*") {
        $r = $p.Range
        # Insert a real paragraph break plus the second line of the address
        # right after the existing run; this leaves the paragraph that
        # already follows the address (an empty spacer paragraph) untouched.
        $r.InsertAfter("`rSanta Clara, CA 95051")
        break
    }
}

# Trim the original paragraph down to just the street address.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*2900 Sanor Pl, Santa Clara CA 95051*") {
        $p.Range.Text = "2900 Sanor Pl"
        break
    }
}

# The newly-created "Santa Clara, CA 95051" paragraph was synthesized from
# scratch, so Word doesn't mark its <w:t> with xml:space="preserve" the way
# the rest of this document's runs are marked. Re-author that single
# paragraph's OOXML directly so it matches the surrounding formatting
# exactly (same rPr/pPr as its sibling runs) including that attribute.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Santa Clara, CA 95051*") {
        $r = $p.Range
        $xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Santa Clara, CA 95051</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $r.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# 3. Remove the now-unwanted empty "No Spacing" paragraph that used to
#    follow "Board of Directors".
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Board of Directors*") {
        $nxt = $p.Next()
        if ($nxt -ne $null) {
            $nxt.Range.Delete()
        }
        break
    }
}
